# Automatic tracker update: fill in results for rows 107 and 109,
# and append two new rows (116 and 117) with new match data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously empty results for rows 107 and 109
$ws.Range("G107").Value = "Fallo"
$ws.Range("H107").Value = -1

$ws.Range("G109").Value = "Fallo"
$ws.Range("H109").Value = -1

# Append new row 116
$ws.Range("A116").Value = 14677775
$ws.Range("B116").NumberFormat = "@"
$ws.Range("B116").Value = "2025-09-13"
$ws.Range("B116").Style = "Normal"
$ws.Range("C116").Value = "Mateus Alves"
$ws.Range("D116").Value = "Bruno Kuzuhara"
$ws.Range("E116").Value = "Gana Bruno Kuzuhara"
$ws.Range("F116").Value = 2.25
$ws.Range("G116").Value = "'"
$ws.Range("G116").Style = "Normal"
$ws.Range("H116").Value = "'"
$ws.Range("H116").Style = "Normal"

# Append new row 117
$ws.Range("A117").Value = 14677776
$ws.Range("B117").NumberFormat = "@"
$ws.Range("B117").Value = "2025-09-13"
$ws.Range("B117").Style = "Normal"
$ws.Range("C117").Value = "Igor Marcondes"
$ws.Range("D117").Value = "Pedro Sakamoto"
$ws.Range("E117").Value = "Gana Igor Marcondes"
$ws.Range("F117").Value = 1.67
$ws.Range("G117").Value = "'"
$ws.Range("G117").Style = "Normal"
$ws.Range("H117").Value = "'"
$ws.Range("H117").Style = "Normal"
